# Integration with ELK stack, Reportportal
# Update RUNMANAGER (testcase master) and DATA (run-manager) sheets so that
# the loginlogoutTest case executes with count 2, while the verifyNewUser
# rows toggle their "run" flag, and refresh the remembered cell selections.

$wb = $excel.ActiveWorkbook

$wsRun  = $wb.Worksheets.Item("RUNMANAGER")
$wsData = $wb.Worksheets.Item("DATA")

# --- RUNMANAGER sheet: flip "execute" to yes and "count" to 2 for loginlogoutTest
$wsRun.Range("C2").Value = "yes"
$wsRun.Range("E2").Value = "'2"

# --- DATA sheet: toggle the "run" flag for the verifyNewUser rows
$wsData.Range("B4").Value = "no"
$wsData.Range("B5").Value = "no"
$wsData.Range("B7").Value = "yes"
$wsData.Range("B9").Value = "yes"

# --- Refresh remembered selections on each sheet, ending with DATA active
$wsRun.Activate()
$null = $wsRun.Range("E4").Select()

$wsData.Activate()
$null = $wsData.Range("A9").Select()
